$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.714.46"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "1.634.26"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "212.25"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "0.522"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "23.26"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +1.13%  "
$ws.Range("E10").Value = "  +0.27%  "
$ws.Range("D11").Value = "0.0861"
$ws.Range("E11").Value = "  -3.18%  "
$ws.Range("D12").Value = "1.864.79"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").Value = "1.634.26"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("D15").Value = "0.554"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("D16").Value = "65.23"
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("D17").Value = "27.671.87"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").Value = "230.15"
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "'10.70"
$ws.Range("E22").Value = "  +4.74%  "
$ws.Range("E23").Value = "  +0.92%  "
$ws.Range("E24").Value = "  +2.47%  "
$ws.Range("D25").Value = "'148.80"
$ws.Range("E25").Value = "  -1.43%  "
$ws.Range("E26").Value = "  -1.14%  "
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("E31").Value = "  -1.02%  "
$ws.Range("E32").Value = "  -0.82%  "
$ws.Range("D33").Value = "1.472.08"
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("D34").Value = "3.08"
$ws.Range("E34").Value = "  -1.08%  "
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("E37").Value = "  +5.68%  "
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("E39").Value = "  -1.60%  "
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("E41").Value = "  +0.78%  "
$ws.Range("D42").Value = "67.95"
$ws.Range("E42").Value = "  -1.59%  "
$ws.Range("D43").Value = "2.47"
$ws.Range("E43").Value = "  +0.89%  "
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("E45").Value = "  -4.62%  "
$ws.Range("D46").Value = "1.774.99"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("E47").Value = "  +1.10%  "
$ws.Range("D48").Value = "87.65"
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("E49").Value = "  -1.41%  "
$ws.Range("D50").Value = "0.0992"
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("D51").Value = "7.69"
$ws.Range("E51").Value = "  -1.91%  "
